# Renamed arcs to links
$wb = $excel.ActiveWorkbook

# Rename the "arcs" worksheet to "links"
$ws = $wb.Worksheets.Item("arcs")
$ws.Name = "links"

# The renamed sheet becomes the active tab (it was the selected tab
# after the rename in the authored workbook).
$ws.Activate()
